$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "71.873.60"
$ws.Cells.Item(2, 5).Value = "  -0.18%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.991.66"
$ws.Cells.Item(3, 5).Value = "  -1.07%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.17%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "541.62"
$ws.Cells.Item(5, 5).Value = "  +4.61%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "149.35"
$ws.Cells.Item(6, 5).Value = "  +0.77%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.696"
$ws.Cells.Item(7, 5).Value = "  +12.16%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.05%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.742"
$ws.Cells.Item(9, 5).Value = "  +1.27%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "0.170"
$ws.Cells.Item(10, 5).Value = "  -1.98%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "ShibaInu"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(11, 4).Value = "0.0000322"
$ws.Cells.Item(11, 5).Value = "  -3.00%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Avalanche"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(12, 4).Value = "49.93"
$ws.Cells.Item(12, 5).Value = "  +4.37%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.77%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "4.637.32"
$ws.Cells.Item(14, 5).Value = "  -0.97%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "4.013.59"
$ws.Cells.Item(15, 5).Value = "  -0.64%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "14.03"
$ws.Cells.Item(16, 5).Value = "  -0.24%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "20.37"

# Row 18
$ws.Cells.Item(18, 5).Value = "  -0.14%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -2.45%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "71.742.64"
$ws.Cells.Item(20, 5).Value = "  -0.41%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "427.17"
$ws.Cells.Item(21, 5).Value = "  -2.00%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "97.05"
$ws.Cells.Item(22, 5).Value = "  -0.36%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "3.49"
$ws.Cells.Item(23, 5).Value = "  -0.79%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "4.24"
$ws.Cells.Item(24, 5).Value = "  +6.02%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "14.21"
$ws.Cells.Item(25, 5).Value = "  -2.52%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "11.25"
$ws.Cells.Item(26, 5).Value = "  -5.55%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -4.55%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +18.75%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +1.07%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "36.72"
$ws.Cells.Item(30, 5).Value = "  -0.43%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).Value = "0.131"
$ws.Cells.Item(31, 5).Value = "  +1.93%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "NEARProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(32, 4).Value = "7.23"
$ws.Cells.Item(32, 5).Value = "  +2.50%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -1.29%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(34, 4).Value = "48.50"
$ws.Cells.Item(34, 5).Value = "  +20.12%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Bittensor"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(35, 4).Value = "673.49"
$ws.Cells.Item(35, 5).Value = "  -3.41%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "65.17"
$ws.Cells.Item(36, 5).Value = "  -4.32%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.440"
$ws.Cells.Item(37, 5).Value = "  +0.96%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.150"
$ws.Cells.Item(38, 5).Value = "  -2.04%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "0.0₃0816"
$ws.Cells.Item(39, 5).Value = "  -7.86%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "3.39"
$ws.Cells.Item(40, 5).Value = "  -7.18%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.11%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "3.35"
$ws.Cells.Item(42, 5).Value = "  +5.73%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.25%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.03%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +3.16%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "9.82"
$ws.Cells.Item(46, 5).Value = "  +8.86%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "2.65"
$ws.Cells.Item(47, 5).Value = "  -3.22%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -5.10%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -3.83%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "0.000269"
$ws.Cells.Item(50, 5).Value = "  +0.36%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "144.00"
$ws.Cells.Item(51, 5).Value = "  +0.91%  "
